$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark rows 108 and 110 as Cancelled (Column D)
$ws.Range("D108").Value = 1
$ws.Range("D110").Value = 1

# Set Weapon for row 112 to "E"
$ws.Range("C112").Value = "E"

# Add new row 113 - "OT" / "F" entry (another copy with different weapon)
$ws.Range("A113").Value = 20211128
$ws.Range("B113").Value = "OT"
$ws.Range("C113").Value = "F"
$ws.Range("D113").Value = 0
$ws.Range("E113").Value = 0
$ws.Range("F113").Value = "13:00"
$ws.Range("G113").Value = "12:30"
$ws.Range("H113").Value = "FSA"

# Match formatting/styles of the row above
$ws.Range("B112:I112").Copy()
$ws.Range("B113:I113").PasteSpecial(-4122)

$ws.Range("K113").Formula = "=IF(Table1[[#This Row],[Cancelled]]=1,""N/A"",Table1[[#This Row],[Date]]&Table1[[#This Row],[Category]]&IF(Table1[[#This Row],[SplitGender]]=1,IF(OR(Table1[[#This Row],[Category]]=""U9"",Table1[[#This Row],[Category]]=""U11"",Table1[[#This Row],[Category]]=""U13""),""B"",""M""),"""")&Table1[[#This Row],[Weapon]])"
$ws.Range("L113").Formula = "=IF(Table1[[#This Row],[Cancelled]]=1,"""",IF(Table1[[#This Row],[SplitGender]]=0,""N/A"",Table1[[#This Row],[Date]]&Table1[[#This Row],[Category]]&IF(Table1[[#This Row],[SplitGender]]=1,IF(OR(Table1[[#This Row],[Category]]=""U9"",Table1[[#This Row],[Category]]=""U11"",Table1[[#This Row],[Category]]=""U13""),""G"",""W""),"""")&Table1[[#This Row],[Weapon]]))"

# Resize the XML table (Table1) to include the new row
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:I113"))
